$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Betas")

$ws.Range("B2").Value = -15155.30149999998
$ws.Range("L2").Value = 0.03000000000000264
$ws.Range("H3").Value = 0
$ws.Range("L3").Value = 0.03300000000000036
$ws.Range("Q3").Value = 1.899999999999993
$ws.Range("H4").Value = 0
$ws.Range("L4").Value = 0.03630000000000022
$ws.Range("Q4").Value = 1.9
$ws.Range("H5").Value = 0
$ws.Range("L5").Value = 0.03993000000000002
$ws.Range("Q5").Value = 1.9
$ws.Range("L6").Value = 0.04392300000000127
$ws.Range("Q6").Value = 1.9
$ws.Range("L7").Value = 0.04831530000000017
$ws.Range("Q7").Value = 1.9
$ws.Range("L8").Value = 0.0531468300000002
$ws.Range("Q8").Value = 1.9
$ws.Range("L9").Value = 0.05846151300000058
$ws.Range("Q9").Value = 1.899999999999999
$ws.Range("L10").Value = 0.06430766429999979
$ws.Range("Q10").Value = 1.899999999999998
$ws.Range("L11").Value = 1.923579476910001
$ws.Range("Q11").Value = 1.9
$ws.Range("L12").Value = 4.519679884700995
$ws.Range("Q12").Value = 4.4937424601
$ws.Range("Q14").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("Q50").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("Q60").Value = 0
$ws.Range("Q61").Value = 0
$ws.Range("Q62").Value = 0
$ws.Range("Q65").Value = 0
$ws.Range("Q67").Value = 2.593742460100002
$ws.Range("Q70").Value = 0
$ws.Range("Q71").Value = 0
$ws.Range("Q72").Value = 0
$ws.Range("Q74").Value = 0
$ws.Range("Q78").Value = 2.593742460100003
$ws.Range("Q81").Value = 0
$ws.Range("Q89").Value = 2.593742460100003
$ws.Range("Q91").Value = 0
$ws.Range("Q100").Value = 2.593742460100003
